$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.100.60'
$ws.Range('E2').Value = '  -0.15%  '

# Row 3
$ws.Range('D3').Value = '1.830.54'
$ws.Range('E3').Value = '  -0.78%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  +0.12%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.37'
$ws.Range('E5').Value = '  -2.08%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6838'
$ws.Range('E6').Value = '  -1.93%  '

# Row 7
$ws.Range('E7').Value = '  +0.20%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3009'
$ws.Range('E8').Value = '  -1.61%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07446'

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.03'
$ws.Range('E10').Value = '  -2.29%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07658'
$ws.Range('E11').Value = '  -2.09%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.840.03'
$ws.Range('E12').Value = '  -0.16%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.048'
$ws.Range('E13').Value = '  -1.41%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6798'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '87.43'
$ws.Range('E15').Value = '  -5.89%  '

# Row 16
$ws.Range('E16').Value = '  -7.61%  '

# Row 17
$ws.Range('D17').Value = '29.102.76'
$ws.Range('E17').Value = '  -0.05%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008170'
$ws.Range('E18').Value = '  -1.46%  '

# Row 19
$ws.Range('D19').Value = '2.081.05'
$ws.Range('E19').Value = '  +0.08%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '227.70'
$ws.Range('E20').Value = '  -5.84%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.49'

# Row 22
$ws.Range('E22').Value = '  +0.15%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.407'
$ws.Range('E23').Value = '  -0.98%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.16%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1453'
$ws.Range('E25').Value = '  -3.64%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.15'
$ws.Range('E26').Value = '  +0.62%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.724'
$ws.Range('E27').Value = '  -1.02%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.10'

# Row 29
$ws.Range('E29').Value = '  -2.17%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.275'
$ws.Range('E30').Value = '  +1.08%  '

# Row 32
$ws.Range('E32').Value = '  -0.34%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05157'

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7656'
$ws.Range('E34').Value = '  -2.49%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.837'
$ws.Range('E35').Value = '  -1.38%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.130'
$ws.Range('E36').Value = '  -1.40%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.676'
$ws.Range('E37').Value = '  -0.56%  '

# Row 38
$ws.Range('D38').Value = '1.307.47'
$ws.Range('E38').Value = '  -0.47%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01835'
$ws.Range('E39').Value = '  -1.72%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.711'
$ws.Range('E40').Value = '  +0.19%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9313'
$ws.Range('E41').Value = '  -1.96%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.806'
$ws.Range('E42').Value = '  -4.49%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '104.39'
$ws.Range('E43').Value = '  -3.07%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9994'
$ws.Range('E44').Value = '  +0.07%  '

# Row 45
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000123'
$ws.Range('E45').Value = '  +0.26%  '

# Row 46
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.981.49'
$ws.Range('E46').Value = '  +0.03%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.03'
$ws.Range('E47').Value = '  +1.36%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5196'
$ws.Range('E48').Value = '  +0.38%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.556'
$ws.Range('E49').Value = '  -1.72%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.769'
$ws.Range('E50').Value = '  +0.48%  '

# Row 51
$ws.Range('E51').Value = '  +0.66%  '
